# feat: add 2022-Q3 data
#
# Before:  总计 | 2022-Q2 | 2022-Q1
# After:   总计 | 2022-Q3 | 2022-Q2 | 2022-Q1
#
# - A new "2022-Q3" sheet is inserted right after "总计" (i.e. before the
#   existing "2022-Q2" sheet), carrying fund-holding data for the new quarter.
# - The "总计" (summary) sheet gets a new top data row for 2022-Q3, with the
#   existing 2022-Q2 / 2022-Q1 rows pushed down (their "序号" index in column
#   A is renumbered accordingly).
# - "2022-Q2" and "2022-Q1" sheets themselves are left content-unchanged,
#   they simply shift one position to the right in the tab strip.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write `text` into $cellRef as a genuine text cell (not a number),
# even when `text` looks numeric (e.g. "3.47", "014283"). Direct
# `Range.Value = "3.47"` assignment gets silently coerced to a number by
# this engine (matching real Excel's "smart" literal entry), so instead we
# stage the string in a scratch cell forced to text via NumberFormat "@",
# then copy/paste-special just the *value* across - which carries the text
# type without dragging the scratch cell's number-format style along.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($ws, $cellRef, $text)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# =======================================================================
# 1) "总计" sheet: insert the 2022-Q3 row at the top of the data, shifting
#    2022-Q2 -> row 3 and 2022-Q1 -> row 4, renumbering column A.
# =======================================================================
$wsTotal = $wb.Worksheets.Item(1)

# Work bottom-up so we never overwrite a source row before it's copied down.
# Row4 (new) <- old Row3 (2022-Q1), with its index renumbered to 2.
$wsTotal.Range("B4").Value = $wsTotal.Range("B3").Value2
$wsTotal.Range("C4").Value = $wsTotal.Range("C3").Value2
$wsTotal.Range("D4").Value = $wsTotal.Range("D3").Value2
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)   # xlPasteFormats (keep the A-column style)
$wsTotal.Range("A4").Value = 2

# Row3 (new) <- old Row2 (2022-Q2), with its index renumbered to 1.
$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value2
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value2
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value2
$wsTotal.Range("A3").Value = 1

# Row2 (new) <- the new 2022-Q3 summary figures, index 0 (A2 already has the
# right style/value from the original row).
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.55

# =======================================================================
# 2) Insert the "2022-Q3" worksheet. Duplicating the existing "2022-Q2"
#    sheet (same columns/header/styles) and placing the copy right before
#    it is the most reliable way to match its formatting exactly.
# =======================================================================
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# 2022-Q2 only has 3 data rows (rows 2-4); 2022-Q3 needs 5 (rows 2-6), so
# stamp the formatting of the last data row onto the two extra rows first.
$wsQ3.Range("A4:H4").Copy()
$wsQ3.Range("A5:H5").PasteSpecial(-4122)   # xlPasteFormats
$wsQ3.Range("A4:H4").Copy()
$wsQ3.Range("A6:H6").PasteSpecial(-4122)   # xlPasteFormats

# ---- row 2 ----
$wsQ3.Range("A2").Value = 0
Set-TextValue $wsQ3 "B2" "014283"
Set-TextValue $wsQ3 "C2" "华夏北交所创新中小企业精选两年定开混合"
Set-TextValue $wsQ3 "D2" "3.47"
Set-TextValue $wsQ3 "E2" "71.59"
Set-TextValue $wsQ3 "F2" "5.69"
Set-TextValue $wsQ3 "G2" "0.1974"
$wsQ3.Range("H2").Value = 5

# ---- row 3 ----
$wsQ3.Range("A3").Value = 1
Set-TextValue $wsQ3 "B3" "014269"
Set-TextValue $wsQ3 "C3" "嘉实北交所精选两年定期混合A"
Set-TextValue $wsQ3 "D3" "2.72"
Set-TextValue $wsQ3 "E3" "90.37"
Set-TextValue $wsQ3 "F3" "6.17"
Set-TextValue $wsQ3 "G3" "0.1678"
$wsQ3.Range("H3").Value = 6

# ---- row 4 ----
$wsQ3.Range("A4").Value = 2
Set-TextValue $wsQ3 "B4" "014279"
Set-TextValue $wsQ3 "C4" "汇添富北交所创新精选两年定开混合A"
Set-TextValue $wsQ3 "D4" "3.20"
Set-TextValue $wsQ3 "E4" "93.27"
Set-TextValue $wsQ3 "F4" "3.97"
Set-TextValue $wsQ3 "G4" "0.1270"
$wsQ3.Range("H4").Value = 9

# ---- row 5 ----
$wsQ3.Range("A5").Value = 3
Set-TextValue $wsQ3 "B5" "014270"
Set-TextValue $wsQ3 "C5" "嘉实北交所精选两年定期混合C"
Set-TextValue $wsQ3 "D5" "0.53"
Set-TextValue $wsQ3 "E5" "90.37"
Set-TextValue $wsQ3 "F5" "6.17"
Set-TextValue $wsQ3 "G5" "0.0327"
$wsQ3.Range("H5").Value = 6

# ---- row 6 ----
$wsQ3.Range("A6").Value = 4
Set-TextValue $wsQ3 "B6" "014280"
Set-TextValue $wsQ3 "C6" "汇添富北交所创新精选两年定开混合C"
Set-TextValue $wsQ3 "D6" "0.51"
Set-TextValue $wsQ3 "E6" "93.27"
Set-TextValue $wsQ3 "F6" "3.97"
Set-TextValue $wsQ3 "G6" "0.0202"
$wsQ3.Range("H6").Value = 9
